$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.215.63'
$ws.Range("E2").Value = '  -1.90%  '

# Row 3
$ws.Range("D3").Value = '3.014.48'
$ws.Range("E3").Value = '  -4.67%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''571.57'
$ws.Range("E5").Value = '  -2.43%  '

# Row 6
$ws.Range("D6").Value = '''128.95'
$ws.Range("E6").Value = '  -4.34%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").Value = '3.018.77'
$ws.Range("E8").Value = '  -4.46%  '

# Row 9
$ws.Range("D9").Value = '''0.497'
$ws.Range("E9").Value = '  -1.89%  '

# Row 10
$ws.Range("E10").Value = '  -3.68%  '

# Row 11
$ws.Range("D11").Value = '''5.20'
$ws.Range("E11").Value = '  -1.02%  '

# Row 12
$ws.Range("D12").Value = '''0.432'
$ws.Range("E12").Value = '  -4.78%  '

# Row 13
$ws.Range("D13").Value = '''0.0000228'
$ws.Range("E13").Value = '  -2.41%  '

# Row 14
$ws.Range("D14").Value = '''33.28'
$ws.Range("E14").Value = '  +0.06%  '

# Row 15
$ws.Range("D15").Value = '''0.119'
$ws.Range("E15").Value = '  +0.41%  '

# Row 16
$ws.Range("D16").Value = '3.519.89'
$ws.Range("E16").Value = '  -4.50%  '

# Row 17
$ws.Range("D17").Value = '61.360.98'
$ws.Range("E17").Value = '  -1.65%  '

# Row 18
$ws.Range("D18").Value = '3.024.65'
$ws.Range("E18").Value = '  -4.58%  '

# Row 19
$ws.Range("D19").Value = '''6.29'
$ws.Range("E19").Value = '  -3.47%  '

# Row 20
$ws.Range("D20").Value = '''439.27'
$ws.Range("E20").Value = '  -3.28%  '

# Row 21
$ws.Range("D21").Value = '''13.25'
$ws.Range("E21").Value = '  -4.77%  '

# Row 22
$ws.Range("D22").Value = '''0.664'
$ws.Range("E22").Value = '  -5.04%  '

# Row 23
$ws.Range("D23").Value = '''7.21'
$ws.Range("E23").Value = '  -5.02%  '

# Row 24
$ws.Range("D24").Value = '''79.95'
$ws.Range("E24").Value = '  -4.23%  '

# Row 25
$ws.Range("D25").Value = '''12.64'
$ws.Range("E25").Value = '  -4.43%  '

# Row 26
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.07%  '

# Row 27
$ws.Range("E27").Value = '  -0.13%  '

# Row 28
$ws.Range("D28").Value = '''2.52'
$ws.Range("E28").Value = '  -5.93%  '

# Row 29
$ws.Range("D29").Value = '''1.97'
$ws.Range("E29").Value = '  -1.84%  '

# Row 30
$ws.Range("D30").Value = '''7.36'
$ws.Range("E30").Value = '  -4.51%  '

# Row 31
$ws.Range("D31").Value = '''6.28'
$ws.Range("E31").Value = '  -7.95%  '

# Row 32
$ws.Range("D32").Value = '''25.58'
$ws.Range("E32").Value = '  -5.67%  '

# Row 33
$ws.Range("D33").Value = '''0.0952'
$ws.Range("E33").Value = '  -7.57%  '

# Row 34
$ws.Range("D34").Value = '''2.29'
$ws.Range("E34").Value = '  -3.60%  '

# Row 35
$ws.Range("D35").Value = '''0.963'
$ws.Range("E35").Value = '  -6.55%  '

# Row 36
$ws.Range("D36").Value = '''5.61'
$ws.Range("E36").Value = '  -5.05%  '

# Row 37
$ws.Range("D37").Value = '''50.19'
$ws.Range("E37").Value = '  -1.80%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0682'
$ws.Range("E38").Value = '  -1.62%  '

# Row 39
$ws.Range("D39").Value = '''0.0368'
$ws.Range("E39").Value = '  -3.75%  '

# Row 40
$ws.Range("D40").Value = '''7.81'
$ws.Range("E40").Value = '  -2.12%  '

# Row 41
$ws.Range("D41").Value = '''0.109'
$ws.Range("E41").Value = '  -2.20%  '

# Row 42
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '''375.04'
$ws.Range("E42").Value = '  -5.11%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '''2.50'
$ws.Range("E43").Value = '  -8.08%  '

# Row 44
$ws.Range("D44").Value = '2.657.86'
$ws.Range("E44").Value = '  -4.89%  '

# Row 45
$ws.Range("E45").Value = '  -0.03%  '

# Row 46
$ws.Range("D46").Value = '''121.91'
$ws.Range("E46").Value = '  -2.55%  '

# Row 47
$ws.Range("E47").Value = '  -4.99%  '

# Row 48
$ws.Range("D48").Value = '''33.76'
$ws.Range("E48").Value = '  -4.81%  '

# Row 49
$ws.Range("D49").Value = '''1.97'
$ws.Range("E49").Value = '  -7.14%  '

# Row 50
$ws.Range("E50").Value = '  -3.01%  '

# Row 51
$ws.Range("D51").Value = '''23.50'
$ws.Range("E51").Value = '  -6.85%  '
